# Applies the weekly Fruta/Hortalizas refresh for "Vega Modelo de Temuco - Arandano (blue)":
# - rows 29-72 keep their row position but get refreshed Fecha/Calidad/Volumen/Precios/Unidad/Origen data
# - two brand-new data rows (73, 74) are appended at the bottom
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 29
$ws.Cells.Item(29, 4).Value = 44540
$ws.Cells.Item(29, 12).Value = 'Primera'
$ws.Cells.Item(29, 13).Value = 200
$ws.Cells.Item(29, 14).Value = 2800
$ws.Cells.Item(29, 15).Value = 2800
$ws.Cells.Item(29, 16).Value = 2800
$ws.Cells.Item(29, 17).Value = '$/kilo'
$ws.Cells.Item(29, 18).Value = 'Región del Maule'
$ws.Cells.Item(29, 19).Value = 2800
$ws.Cells.Item(29, 20).Value = 1

# Row 30
$ws.Cells.Item(30, 4).Value = 44540
$ws.Cells.Item(30, 12).Value = 'Segunda'
$ws.Cells.Item(30, 13).Value = 180
$ws.Cells.Item(30, 14).Value = 2500
$ws.Cells.Item(30, 15).Value = 2500
$ws.Cells.Item(30, 16).Value = 2500
$ws.Cells.Item(30, 17).Value = '$/kilo'
$ws.Cells.Item(30, 18).Value = 'Región del Maule'
$ws.Cells.Item(30, 19).Value = 2500
$ws.Cells.Item(30, 20).Value = 1

# Row 31
$ws.Cells.Item(31, 4).Value = 44246
$ws.Cells.Item(31, 12).Value = 'Primera'
$ws.Cells.Item(31, 13).Value = 150
$ws.Cells.Item(31, 14).Value = 1500
$ws.Cells.Item(31, 15).Value = 1500
$ws.Cells.Item(31, 16).Value = 1500
$ws.Cells.Item(31, 17).Value = '$/kilo'
$ws.Cells.Item(31, 18).Value = 'Región de La Araucanía'
$ws.Cells.Item(31, 19).Value = 1500
$ws.Cells.Item(31, 20).Value = 1

# Row 32
$ws.Cells.Item(32, 4).Value = 44249
$ws.Cells.Item(32, 12).Value = 'Primera'
$ws.Cells.Item(32, 13).Value = 200
$ws.Cells.Item(32, 14).Value = 1500
$ws.Cells.Item(32, 15).Value = 1500
$ws.Cells.Item(32, 16).Value = 1500
$ws.Cells.Item(32, 17).Value = '$/kilo'
$ws.Cells.Item(32, 18).Value = 'Región de La Araucanía'
$ws.Cells.Item(32, 19).Value = 1500
$ws.Cells.Item(32, 20).Value = 1

# Row 33
$ws.Cells.Item(33, 4).Value = 44172
$ws.Cells.Item(33, 12).Value = 'Primera'
$ws.Cells.Item(33, 13).Value = 300
$ws.Cells.Item(33, 14).Value = 1800
$ws.Cells.Item(33, 15).Value = 2000
$ws.Cells.Item(33, 16).Value = 1920
$ws.Cells.Item(33, 17).Value = '$/kilo'
$ws.Cells.Item(33, 18).Value = 'Región de La Araucanía'
$ws.Cells.Item(33, 19).Value = 1920
$ws.Cells.Item(33, 20).Value = 1

# Row 34
$ws.Cells.Item(34, 4).Value = 44172
$ws.Cells.Item(34, 12).Value = 'Segunda'
$ws.Cells.Item(34, 13).Value = 200
$ws.Cells.Item(34, 14).Value = 1500
$ws.Cells.Item(34, 15).Value = 1500
$ws.Cells.Item(34, 16).Value = 1500
$ws.Cells.Item(34, 17).Value = '$/kilo'
$ws.Cells.Item(34, 18).Value = 'Región del Maule'
$ws.Cells.Item(34, 19).Value = 1500
$ws.Cells.Item(34, 20).Value = 1

# Row 35
$ws.Cells.Item(35, 4).Value = 44186
$ws.Cells.Item(35, 12).Value = 'Primera'
$ws.Cells.Item(35, 13).Value = 200
$ws.Cells.Item(35, 14).Value = 1500
$ws.Cells.Item(35, 15).Value = 1500
$ws.Cells.Item(35, 16).Value = 1500
$ws.Cells.Item(35, 17).Value = '$/kilo'
$ws.Cells.Item(35, 18).Value = 'Región del Maule'
$ws.Cells.Item(35, 19).Value = 1500
$ws.Cells.Item(35, 20).Value = 1

# Row 36
$ws.Cells.Item(36, 4).Value = 44186
$ws.Cells.Item(36, 12).Value = 'Segunda'
$ws.Cells.Item(36, 13).Value = 250
$ws.Cells.Item(36, 14).Value = 1300
$ws.Cells.Item(36, 15).Value = 1300
$ws.Cells.Item(36, 16).Value = 1300
$ws.Cells.Item(36, 17).Value = '$/kilo'
$ws.Cells.Item(36, 18).Value = 'Región del Maule'
$ws.Cells.Item(36, 19).Value = 1300
$ws.Cells.Item(36, 20).Value = 1

# Row 37
$ws.Cells.Item(37, 4).Value = 44179
$ws.Cells.Item(37, 12).Value = 'Primera'
$ws.Cells.Item(37, 13).Value = 500
$ws.Cells.Item(37, 14).Value = 1600
$ws.Cells.Item(37, 15).Value = 1800
$ws.Cells.Item(37, 16).Value = 1680
$ws.Cells.Item(37, 17).Value = '$/kilo'
$ws.Cells.Item(37, 18).Value = 'Región del Maule'
$ws.Cells.Item(37, 19).Value = 1680
$ws.Cells.Item(37, 20).Value = 1

# Row 38
$ws.Cells.Item(38, 4).Value = 44216
$ws.Cells.Item(38, 12).Value = 'Primera'
$ws.Cells.Item(38, 13).Value = 430
$ws.Cells.Item(38, 14).Value = 1400
$ws.Cells.Item(38, 15).Value = 1500
$ws.Cells.Item(38, 16).Value = 1458
$ws.Cells.Item(38, 17).Value = '$/kilo'
$ws.Cells.Item(38, 18).Value = 'Región del Maule'
$ws.Cells.Item(38, 19).Value = 1458
$ws.Cells.Item(38, 20).Value = 1

# Row 39
$ws.Cells.Item(39, 4).Value = 44539
$ws.Cells.Item(39, 12).Value = 'Primera'
$ws.Cells.Item(39, 13).Value = 350
$ws.Cells.Item(39, 14).Value = 2800
$ws.Cells.Item(39, 15).Value = 3000
$ws.Cells.Item(39, 16).Value = 2886
$ws.Cells.Item(39, 17).Value = '$/kilo'
$ws.Cells.Item(39, 18).Value = 'Región del Maule'
$ws.Cells.Item(39, 19).Value = 2886
$ws.Cells.Item(39, 20).Value = 1

# Row 40
$ws.Cells.Item(40, 4).Value = 44539
$ws.Cells.Item(40, 12).Value = 'Segunda'
$ws.Cells.Item(40, 13).Value = 300
$ws.Cells.Item(40, 14).Value = 2300
$ws.Cells.Item(40, 15).Value = 2300
$ws.Cells.Item(40, 16).Value = 2300
$ws.Cells.Item(40, 17).Value = '$/kilo'
$ws.Cells.Item(40, 18).Value = 'Región del Maule'
$ws.Cells.Item(40, 19).Value = 2300
$ws.Cells.Item(40, 20).Value = 1

# Row 41
$ws.Cells.Item(41, 4).Value = 44176
$ws.Cells.Item(41, 12).Value = 'Primera'
$ws.Cells.Item(41, 13).Value = 80
$ws.Cells.Item(41, 14).Value = 2000
$ws.Cells.Item(41, 15).Value = 2000
$ws.Cells.Item(41, 16).Value = 2000
$ws.Cells.Item(41, 17).Value = '$/kilo'
$ws.Cells.Item(41, 18).Value = 'Región del Maule'
$ws.Cells.Item(41, 19).Value = 2000
$ws.Cells.Item(41, 20).Value = 1

# Row 42
$ws.Cells.Item(42, 4).Value = 44176
$ws.Cells.Item(42, 12).Value = 'Segunda'
$ws.Cells.Item(42, 13).Value = 200
$ws.Cells.Item(42, 14).Value = 1500
$ws.Cells.Item(42, 15).Value = 1500
$ws.Cells.Item(42, 16).Value = 1500
$ws.Cells.Item(42, 17).Value = '$/kilo'
$ws.Cells.Item(42, 18).Value = 'Región del Maule'
$ws.Cells.Item(42, 19).Value = 1500
$ws.Cells.Item(42, 20).Value = 1

# Row 43
$ws.Cells.Item(43, 4).Value = 44165
$ws.Cells.Item(43, 12).Value = 'Primera'
$ws.Cells.Item(43, 13).Value = 300
$ws.Cells.Item(43, 14).Value = 2000
$ws.Cells.Item(43, 15).Value = 2500
$ws.Cells.Item(43, 16).Value = 2167
$ws.Cells.Item(43, 17).Value = '$/kilo'
$ws.Cells.Item(43, 18).Value = 'Región del Maule'
$ws.Cells.Item(43, 19).Value = 2167
$ws.Cells.Item(43, 20).Value = 1

# Row 44
$ws.Cells.Item(44, 4).Value = 44530
$ws.Cells.Item(44, 12).Value = 'Primera'
$ws.Cells.Item(44, 13).Value = 200
$ws.Cells.Item(44, 14).Value = 3000
$ws.Cells.Item(44, 15).Value = 3000
$ws.Cells.Item(44, 16).Value = 3000
$ws.Cells.Item(44, 17).Value = '$/kilo'
$ws.Cells.Item(44, 18).Value = 'Región del Maule'
$ws.Cells.Item(44, 19).Value = 3000
$ws.Cells.Item(44, 20).Value = 1

# Row 45
$ws.Cells.Item(45, 4).Value = 44252
$ws.Cells.Item(45, 12).Value = 'Primera'
$ws.Cells.Item(45, 13).Value = 550
$ws.Cells.Item(45, 14).Value = 1400
$ws.Cells.Item(45, 15).Value = 1500
$ws.Cells.Item(45, 16).Value = 1473
$ws.Cells.Item(45, 17).Value = '$/kilo'
$ws.Cells.Item(45, 18).Value = 'Región de La Araucanía'
$ws.Cells.Item(45, 19).Value = 1473
$ws.Cells.Item(45, 20).Value = 1

# Row 46
$ws.Cells.Item(46, 4).Value = 44245
$ws.Cells.Item(46, 12).Value = 'Primera'
$ws.Cells.Item(46, 13).Value = 200
$ws.Cells.Item(46, 14).Value = 1500
$ws.Cells.Item(46, 15).Value = 1500
$ws.Cells.Item(46, 16).Value = 1500
$ws.Cells.Item(46, 17).Value = '$/kilo'
$ws.Cells.Item(46, 18).Value = 'Región de La Araucanía'
$ws.Cells.Item(46, 19).Value = 1500
$ws.Cells.Item(46, 20).Value = 1

# Row 47
$ws.Cells.Item(47, 4).Value = 44243
$ws.Cells.Item(47, 12).Value = 'Primera'
$ws.Cells.Item(47, 13).Value = 150
$ws.Cells.Item(47, 14).Value = 1500
$ws.Cells.Item(47, 15).Value = 1500
$ws.Cells.Item(47, 16).Value = 1500
$ws.Cells.Item(47, 17).Value = '$/kilo'
$ws.Cells.Item(47, 18).Value = 'Región de La Araucanía'
$ws.Cells.Item(47, 19).Value = 1500
$ws.Cells.Item(47, 20).Value = 1

# Row 48
$ws.Cells.Item(48, 4).Value = 44250
$ws.Cells.Item(48, 12).Value = 'Primera'
$ws.Cells.Item(48, 13).Value = 60
$ws.Cells.Item(48, 14).Value = 1500
$ws.Cells.Item(48, 15).Value = 1500
$ws.Cells.Item(48, 16).Value = 1500
$ws.Cells.Item(48, 17).Value = '$/kilo'
$ws.Cells.Item(48, 18).Value = 'Región de La Araucanía'
$ws.Cells.Item(48, 19).Value = 1500
$ws.Cells.Item(48, 20).Value = 1

# Row 49
$ws.Cells.Item(49, 4).Value = 44503
$ws.Cells.Item(49, 12).Value = 'Primera'
$ws.Cells.Item(49, 13).Value = 300
$ws.Cells.Item(49, 14).Value = 4500
$ws.Cells.Item(49, 15).Value = 4500
$ws.Cells.Item(49, 16).Value = 4500
$ws.Cells.Item(49, 17).Value = '$/kilo'
$ws.Cells.Item(49, 18).Value = 'Región de La Araucanía'
$ws.Cells.Item(49, 19).Value = 4500
$ws.Cells.Item(49, 20).Value = 1

# Row 50
$ws.Cells.Item(50, 4).Value = 44188
$ws.Cells.Item(50, 12).Value = 'Primera'
$ws.Cells.Item(50, 13).Value = 250
$ws.Cells.Item(50, 14).Value = 1500
$ws.Cells.Item(50, 15).Value = 1500
$ws.Cells.Item(50, 16).Value = 1500
$ws.Cells.Item(50, 17).Value = '$/kilo'
$ws.Cells.Item(50, 18).Value = 'Región del Maule'
$ws.Cells.Item(50, 19).Value = 1500
$ws.Cells.Item(50, 20).Value = 1

# Row 51
$ws.Cells.Item(51, 4).Value = 44159
$ws.Cells.Item(51, 12).Value = 'Primera'
$ws.Cells.Item(51, 13).Value = 250
$ws.Cells.Item(51, 14).Value = 3000
$ws.Cells.Item(51, 15).Value = 3000
$ws.Cells.Item(51, 16).Value = 3000
$ws.Cells.Item(51, 17).Value = '$/kilo'
$ws.Cells.Item(51, 18).Value = 'Región del Maule'
$ws.Cells.Item(51, 19).Value = 3000
$ws.Cells.Item(51, 20).Value = 1

# Row 52
$ws.Cells.Item(52, 4).Value = 44159
$ws.Cells.Item(52, 12).Value = 'Segunda'
$ws.Cells.Item(52, 13).Value = 150
$ws.Cells.Item(52, 14).Value = 2800
$ws.Cells.Item(52, 15).Value = 2800
$ws.Cells.Item(52, 16).Value = 2800
$ws.Cells.Item(52, 17).Value = '$/kilo'
$ws.Cells.Item(52, 18).Value = 'Región del Maule'
$ws.Cells.Item(52, 19).Value = 2800
$ws.Cells.Item(52, 20).Value = 1

# Row 53
$ws.Cells.Item(53, 4).Value = 44237
$ws.Cells.Item(53, 12).Value = 'Primera'
$ws.Cells.Item(53, 13).Value = 200
$ws.Cells.Item(53, 14).Value = 1400
$ws.Cells.Item(53, 15).Value = 1400
$ws.Cells.Item(53, 16).Value = 1400
$ws.Cells.Item(53, 17).Value = '$/kilo'
$ws.Cells.Item(53, 18).Value = 'Región de La Araucanía'
$ws.Cells.Item(53, 19).Value = 1400
$ws.Cells.Item(53, 20).Value = 1

# Row 54
$ws.Cells.Item(54, 4).Value = 44237
$ws.Cells.Item(54, 12).Value = 'Primera'
$ws.Cells.Item(54, 13).Value = 1000
$ws.Cells.Item(54, 14).Value = 1400
$ws.Cells.Item(54, 15).Value = 1500
$ws.Cells.Item(54, 16).Value = 1450
$ws.Cells.Item(54, 17).Value = '$/kilo'
$ws.Cells.Item(54, 18).Value = 'Región del Maule'
$ws.Cells.Item(54, 19).Value = 1450
$ws.Cells.Item(54, 20).Value = 1

# Row 55
$ws.Cells.Item(55, 4).Value = 44529
$ws.Cells.Item(55, 12).Value = 'Primera'
$ws.Cells.Item(55, 13).Value = 1050
$ws.Cells.Item(55, 14).Value = 3000
$ws.Cells.Item(55, 15).Value = 3200
$ws.Cells.Item(55, 16).Value = 3086
$ws.Cells.Item(55, 17).Value = '$/kilo'
$ws.Cells.Item(55, 18).Value = 'Región del Maule'
$ws.Cells.Item(55, 19).Value = 3086
$ws.Cells.Item(55, 20).Value = 1

# Row 56
$ws.Cells.Item(56, 4).Value = 44505
$ws.Cells.Item(56, 12).Value = 'Primera'
$ws.Cells.Item(56, 13).Value = 30
$ws.Cells.Item(56, 14).Value = 10000
$ws.Cells.Item(56, 15).Value = 10000
$ws.Cells.Item(56, 16).Value = 10000
$ws.Cells.Item(56, 17).Value = '$/bandeja 2 kilos'
$ws.Cells.Item(56, 18).Value = 'Provincia de Limarí'
$ws.Cells.Item(56, 19).Value = 5000
$ws.Cells.Item(56, 20).Value = 2

# Row 57
$ws.Cells.Item(57, 4).Value = 44208
$ws.Cells.Item(57, 12).Value = 'Primera'
$ws.Cells.Item(57, 13).Value = 200
$ws.Cells.Item(57, 14).Value = 1400
$ws.Cells.Item(57, 15).Value = 1400
$ws.Cells.Item(57, 16).Value = 1400
$ws.Cells.Item(57, 17).Value = '$/kilo'
$ws.Cells.Item(57, 18).Value = 'Región del Maule'
$ws.Cells.Item(57, 19).Value = 1400
$ws.Cells.Item(57, 20).Value = 1

# Row 58
$ws.Cells.Item(58, 4).Value = 44210
$ws.Cells.Item(58, 12).Value = 'Primera'
$ws.Cells.Item(58, 13).Value = 800
$ws.Cells.Item(58, 14).Value = 1400
$ws.Cells.Item(58, 15).Value = 1400
$ws.Cells.Item(58, 16).Value = 1400
$ws.Cells.Item(58, 17).Value = '$/kilo'
$ws.Cells.Item(58, 18).Value = 'Región del Maule'
$ws.Cells.Item(58, 19).Value = 1400
$ws.Cells.Item(58, 20).Value = 1

# Row 59
$ws.Cells.Item(59, 4).Value = 44242
$ws.Cells.Item(59, 12).Value = 'Primera'
$ws.Cells.Item(59, 13).Value = 110
$ws.Cells.Item(59, 14).Value = 1500
$ws.Cells.Item(59, 15).Value = 1500
$ws.Cells.Item(59, 16).Value = 1500
$ws.Cells.Item(59, 17).Value = '$/kilo'
$ws.Cells.Item(59, 18).Value = 'Región del Maule'
$ws.Cells.Item(59, 19).Value = 1500
$ws.Cells.Item(59, 20).Value = 1

# Row 60
$ws.Cells.Item(60, 4).Value = 44518
$ws.Cells.Item(60, 12).Value = 'Primera'
$ws.Cells.Item(60, 13).Value = 80
$ws.Cells.Item(60, 14).Value = 3500
$ws.Cells.Item(60, 15).Value = 3500
$ws.Cells.Item(60, 16).Value = 3500
$ws.Cells.Item(60, 17).Value = '$/kilo'
$ws.Cells.Item(60, 18).Value = 'Región del Maule'
$ws.Cells.Item(60, 19).Value = 3500
$ws.Cells.Item(60, 20).Value = 1

# Row 61
$ws.Cells.Item(61, 4).Value = 44168
$ws.Cells.Item(61, 12).Value = 'Primera'
$ws.Cells.Item(61, 13).Value = 700
$ws.Cells.Item(61, 14).Value = 1800
$ws.Cells.Item(61, 15).Value = 2000
$ws.Cells.Item(61, 16).Value = 1886
$ws.Cells.Item(61, 17).Value = '$/kilo'
$ws.Cells.Item(61, 18).Value = 'Región de La Araucanía'
$ws.Cells.Item(61, 19).Value = 1886
$ws.Cells.Item(61, 20).Value = 1

# Row 62
$ws.Cells.Item(62, 4).Value = 44215
$ws.Cells.Item(62, 12).Value = 'Primera'
$ws.Cells.Item(62, 13).Value = 150
$ws.Cells.Item(62, 14).Value = 1400
$ws.Cells.Item(62, 15).Value = 1400
$ws.Cells.Item(62, 16).Value = 1400
$ws.Cells.Item(62, 17).Value = '$/kilo'
$ws.Cells.Item(62, 18).Value = 'Región del Maule'
$ws.Cells.Item(62, 19).Value = 1400
$ws.Cells.Item(62, 20).Value = 1

# Row 63
$ws.Cells.Item(63, 4).Value = 44204
$ws.Cells.Item(63, 12).Value = 'Primera'
$ws.Cells.Item(63, 13).Value = 150
$ws.Cells.Item(63, 14).Value = 1500
$ws.Cells.Item(63, 15).Value = 1500
$ws.Cells.Item(63, 16).Value = 1500
$ws.Cells.Item(63, 17).Value = '$/kilo'
$ws.Cells.Item(63, 18).Value = 'Región del Maule'
$ws.Cells.Item(63, 19).Value = 1500
$ws.Cells.Item(63, 20).Value = 1

# Row 64
$ws.Cells.Item(64, 4).Value = 44244
$ws.Cells.Item(64, 12).Value = 'Primera'
$ws.Cells.Item(64, 13).Value = 200
$ws.Cells.Item(64, 14).Value = 1500
$ws.Cells.Item(64, 15).Value = 1500
$ws.Cells.Item(64, 16).Value = 1500
$ws.Cells.Item(64, 17).Value = '$/kilo'
$ws.Cells.Item(64, 18).Value = 'Región de La Araucanía'
$ws.Cells.Item(64, 19).Value = 1500
$ws.Cells.Item(64, 20).Value = 1

# Row 65
$ws.Cells.Item(65, 4).Value = 44166
$ws.Cells.Item(65, 12).Value = 'Primera'
$ws.Cells.Item(65, 13).Value = 100
$ws.Cells.Item(65, 14).Value = 2500
$ws.Cells.Item(65, 15).Value = 2500
$ws.Cells.Item(65, 16).Value = 2500
$ws.Cells.Item(65, 17).Value = '$/kilo'
$ws.Cells.Item(65, 18).Value = 'Región del Maule'
$ws.Cells.Item(65, 19).Value = 2500
$ws.Cells.Item(65, 20).Value = 1

# Row 66
$ws.Cells.Item(66, 4).Value = 44522
$ws.Cells.Item(66, 12).Value = 'Primera'
$ws.Cells.Item(66, 13).Value = 200
$ws.Cells.Item(66, 14).Value = 3200
$ws.Cells.Item(66, 15).Value = 3200
$ws.Cells.Item(66, 16).Value = 3200
$ws.Cells.Item(66, 17).Value = '$/kilo'
$ws.Cells.Item(66, 18).Value = 'Región de La Araucanía'
$ws.Cells.Item(66, 19).Value = 3200
$ws.Cells.Item(66, 20).Value = 1

# Row 67
$ws.Cells.Item(67, 4).Value = 44218
$ws.Cells.Item(67, 12).Value = 'Primera'
$ws.Cells.Item(67, 13).Value = 250
$ws.Cells.Item(67, 14).Value = 1400
$ws.Cells.Item(67, 15).Value = 1400
$ws.Cells.Item(67, 16).Value = 1400
$ws.Cells.Item(67, 17).Value = '$/kilo'
$ws.Cells.Item(67, 18).Value = 'Región del Maule'
$ws.Cells.Item(67, 19).Value = 1400
$ws.Cells.Item(67, 20).Value = 1

# Row 68
$ws.Cells.Item(68, 4).Value = 44497
$ws.Cells.Item(68, 12).Value = 'Primera'
$ws.Cells.Item(68, 13).Value = 300
$ws.Cells.Item(68, 14).Value = 5000
$ws.Cells.Item(68, 15).Value = 5000
$ws.Cells.Item(68, 16).Value = 5000
$ws.Cells.Item(68, 17).Value = '$/kilo'
$ws.Cells.Item(68, 18).Value = 'Región del Maule'
$ws.Cells.Item(68, 19).Value = 5000
$ws.Cells.Item(68, 20).Value = 1

# Row 69
$ws.Cells.Item(69, 4).Value = 44525
$ws.Cells.Item(69, 12).Value = 'Primera'
$ws.Cells.Item(69, 13).Value = 1450
$ws.Cells.Item(69, 14).Value = 3000
$ws.Cells.Item(69, 15).Value = 3200
$ws.Cells.Item(69, 16).Value = 3110
$ws.Cells.Item(69, 17).Value = '$/kilo'
$ws.Cells.Item(69, 18).Value = 'Región del Maule'
$ws.Cells.Item(69, 19).Value = 3110
$ws.Cells.Item(69, 20).Value = 1

# Row 70
$ws.Cells.Item(70, 4).Value = 44508
$ws.Cells.Item(70, 12).Value = 'Primera'
$ws.Cells.Item(70, 13).Value = 120
$ws.Cells.Item(70, 14).Value = 4000
$ws.Cells.Item(70, 15).Value = 4000
$ws.Cells.Item(70, 16).Value = 4000
$ws.Cells.Item(70, 17).Value = '$/kilo'
$ws.Cells.Item(70, 18).Value = 'Provincia de Limarí'
$ws.Cells.Item(70, 19).Value = 4000
$ws.Cells.Item(70, 20).Value = 1

# Row 71
$ws.Cells.Item(71, 4).Value = 44160
$ws.Cells.Item(71, 12).Value = 'Primera'
$ws.Cells.Item(71, 13).Value = 450
$ws.Cells.Item(71, 14).Value = 2500
$ws.Cells.Item(71, 15).Value = 2500
$ws.Cells.Item(71, 16).Value = 2500
$ws.Cells.Item(71, 17).Value = '$/kilo'
$ws.Cells.Item(71, 18).Value = 'Región del Maule'
$ws.Cells.Item(71, 19).Value = 2500
$ws.Cells.Item(71, 20).Value = 1

# Row 72
$ws.Cells.Item(72, 4).Value = 44211
$ws.Cells.Item(72, 12).Value = 'Primera'
$ws.Cells.Item(72, 13).Value = 200
$ws.Cells.Item(72, 14).Value = 1400
$ws.Cells.Item(72, 15).Value = 1400
$ws.Cells.Item(72, 16).Value = 1400
$ws.Cells.Item(72, 17).Value = '$/kilo'
$ws.Cells.Item(72, 18).Value = 'Región del Maule'
$ws.Cells.Item(72, 19).Value = 1400
$ws.Cells.Item(72, 20).Value = 1

# Row 73
$ws.Cells.Item(73, 1).Value = 10
$ws.Cells.Item(73, 2).Value = 'Vega Modelo de Temuco'
$ws.Cells.Item(73, 3).Value = 'La Araucanía'
$ws.Cells.Item(73, 5).Value = 9
$ws.Cells.Item(73, 6).Value = 'Fruta'
$ws.Cells.Item(73, 7).Value = 100101
$ws.Cells.Item(73, 8).Value = 'Berries'
$ws.Cells.Item(73, 9).Value = 100101001
$ws.Cells.Item(73, 10).Value = 'Arándano (blue)'
$ws.Cells.Item(73, 11).Value = 'Sin especificar'
$ws.Cells.Item(73, 4).Value = 44167
$ws.Cells.Item(73, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(73, 12).Value = 'Primera'
$ws.Cells.Item(73, 13).Value = 200
$ws.Cells.Item(73, 14).Value = 2000
$ws.Cells.Item(73, 15).Value = 2000
$ws.Cells.Item(73, 16).Value = 2000
$ws.Cells.Item(73, 17).Value = '$/kilo'
$ws.Cells.Item(73, 18).Value = 'Región del Maule'
$ws.Cells.Item(73, 19).Value = 2000
$ws.Cells.Item(73, 20).Value = 1

# Row 74
$ws.Cells.Item(74, 1).Value = 10
$ws.Cells.Item(74, 2).Value = 'Vega Modelo de Temuco'
$ws.Cells.Item(74, 3).Value = 'La Araucanía'
$ws.Cells.Item(74, 5).Value = 9
$ws.Cells.Item(74, 6).Value = 'Fruta'
$ws.Cells.Item(74, 7).Value = 100101
$ws.Cells.Item(74, 8).Value = 'Berries'
$ws.Cells.Item(74, 9).Value = 100101001
$ws.Cells.Item(74, 10).Value = 'Arándano (blue)'
$ws.Cells.Item(74, 11).Value = 'Sin especificar'
$ws.Cells.Item(74, 4).Value = 44189
$ws.Cells.Item(74, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(74, 12).Value = 'Primera'
$ws.Cells.Item(74, 13).Value = 200
$ws.Cells.Item(74, 14).Value = 1500
$ws.Cells.Item(74, 15).Value = 1500
$ws.Cells.Item(74, 16).Value = 1500
$ws.Cells.Item(74, 17).Value = '$/kilo'
$ws.Cells.Item(74, 18).Value = 'Región del Maule'
$ws.Cells.Item(74, 19).Value = 1500
$ws.Cells.Item(74, 20).Value = 1
